$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = -4
$ws.Range("F4").Value = -5
$ws.Range("F5").Value = -8
$ws.Range("F6").Value = -2
$ws.Range("F8").Value = -3
$ws.Range("F9").Value = -3
$ws.Range("F10").Value = -1
